$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.955.03'
$ws.Range('E2').Value = '  +1.59%  '
$ws.Range('D3').Value = '1.941.36'
$ws.Range('E3').Value = '  +1.15%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.006'
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '335.61'
$ws.Range('E5').Value = '  +2.86%  '
$ws.Range('E6').Value = '  -0.12%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4843'
$ws.Range('E7').Value = '  +0.32%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4134'
$ws.Range('E8').Value = '  +1.30%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08182'
$ws.Range('E9').Value = '  -0.55%  '
$ws.Range('E10').Value = '  -0.60%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '23.70'
$ws.Range('E11').Value = '  +0.61%  '
$ws.Range('D12').Value = '1.957.17'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.103'
$ws.Range('E13').Value = '  +0.86%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.311'
$ws.Range('E14').Value = '  +1.02%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.22'
$ws.Range('E15').Value = '  -0.16%  '
$ws.Range('E16').Value = '  +0.78%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.008'
$ws.Range('E17').Value = '  -0.17%  '
$ws.Range('E18').Value = '  -0.34%  '
$ws.Range('E19').Value = '  +0.08%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.006'
$ws.Range('E20').Value = '  -0.13%  '
$ws.Range('D21').Value = '29.941.21'
$ws.Range('E21').Value = '  +1.44%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.652'
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.91'
$ws.Range('E23').Value = '  +1.09%  '
$ws.Range('E24').Value = '  -0.45%  '
$ws.Range('D25').Value = '2.172.85'
$ws.Range('E25').Value = '  +0.84%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.660'
$ws.Range('E26').Value = '  +0.33%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '156.74'
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('E28').Value = '  +0.09%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.107'
$ws.Range('E29').Value = '  -0.47%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '121.36'
$ws.Range('E30').Value = '  +0.77%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.012'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09641'
$ws.Range('E33').Value = '  +1.02%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.427'
$ws.Range('E34').Value = '  +3.15%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.548'
$ws.Range('E35').Value = '  -0.43%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06584'
$ws.Range('E36').Value = '  +7.24%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02291'
$ws.Range('E37').Value = '  +0.16%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.214'
$ws.Range('E38').Value = '  +2.57%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5980'
$ws.Range('E39').Value = '  -0.07%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.014'
$ws.Range('E40').Value = '  -0.35%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '10.75'
$ws.Range('E41').Value = '  -0.69%  '
$ws.Range('E42').Value = '  -0.09%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.495'
$ws.Range('E43').Value = '  +3.29%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.272'
$ws.Range('E44').Value = '  -0.64%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '12.34'
$ws.Range('E45').Value = '  -1.00%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.07508'
$ws.Range('E46').Value = '  -1.26%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5583'
$ws.Range('E47').Value = '  -0.06%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.988'
$ws.Range('E48').Value = '  +1.52%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '117.69'
$ws.Range('E49').Value = '  +0.06%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '72.84'
$ws.Range('E50').Value = '  +0.03%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.419'
$ws.Range('E51').Value = '  -0.33%  '
